$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Copy()
$ws.Range("A32").PasteSpecial(-4122)

$ws.Range("A32").Value = 46002
$ws.Range("B32").Value = 96

$ws.Range("A32:B32").Select()

$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
